$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A4: phone number stored as text "79174445" (leading apostrophe forces
# text, matching the source data's inlineStr string type for the phone
# column), then reset the style so it doesn't pick up an extra
# quote-prefix format that the original file doesn't have.
$ws.Range("A4").Formula = "'79174445"
$ws.Range("A4").Style = "Normal"

# B4: empty text value (birthday left blank for this row).
$ws.Range("B4").Formula = "'"
$ws.Range("B4").Style = "Normal"

# C4: points reset to 0.
$ws.Range("C4").Value = 0
